$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Ccl11 -> Ackr2, ECs)
$ws.Range("G2").Value = 0.8417533333333332
$ws.Range("H2").Value = 2.52526
$ws.Range("I2").Value = 0.01079423211523897
$ws.Range("J2").Value = 0.01079423211523897
$ws.Range("Q2").Value = 0.4603453581288889
$ws.Range("R2").Value = 4.14310822316
$ws.Range("S2").Value = 0.01079423211523897
$ws.Range("T2").Value = 0.01079423211523897

# Row 3 (FAPs)
$ws.Range("G3").Value = 69.05064766666666
$ws.Range("I3").Value = 0.885471656726338
$ws.Range("J3").Value = 0.8854716567263378
$ws.Range("Q3").Value = 37.76301663489311
$ws.Range("S3").Value = 0.885471656726338
$ws.Range("T3").Value = 0.8854716567263378

# Row 4 (MuSCs)
$ws.Range("G4").Value = 7.697976666666666
$ws.Range("H4").Value = 23.09393
$ws.Range("I4").Value = 0.09871507918910555
$ws.Range("J4").Value = 0.09871507918910553
$ws.Range("Q4").Value = 4.209936195264445
$ws.Range("R4").Value = 37.88942575738
$ws.Range("S4").Value = 0.09871507918910555
$ws.Range("T4").Value = 0.09871507918910553

# Row 5 (Resolving-Mac)
$ws.Range("G5").Value = 0.391393
$ws.Range("H5").Value = 1.174179
$ws.Range("I5").Value = 0.005019031969317685
$ws.Range("J5").Value = 0.005019031969317684
$ws.Range("Q5").Value = 0.2140483959126667
$ws.Range("R5").Value = 1.926435563214
$ws.Range("S5").Value = 0.005019031969317685
$ws.Range("T5").Value = 0.005019031969317684
